# Auto-generated edit script: updates cryptos price/volume table cells
# to match the refreshed scrape (commit: "Updated cryptos list on Thu Aug 31 09:30:45 UTC 2023 with GitHub Actions").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.335.00"
$ws.Range("E2").Value = "  -0.65%  "
$ws.Range("D3").Value = "1.711.50"
$ws.Range("E3").Value = "  -0.58%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "224.42"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.65%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5291"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.04%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.005"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.10%  "
$ws.Range("B8").Value = "Cardano"
$ws.Range("C8").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2665"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("B9").Value = "Dogecoin"
$ws.Range("C9").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06670"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.31%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.85"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.64%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07669"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.98%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.498"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.54%  "
$ws.Range("D13").Value = "1.947.37"
$ws.Range("E13").Value = "  -0.64%  "
$ws.Range("D14").Value = "1.708.92"
$ws.Range("E14").Value = "  -0.97%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5828"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.00%  "
$ws.Range("D16").Value = "0.0₅8215"
$ws.Range("E16").Value = "  -0.96%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "67.96"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.18%  "
$ws.Range("D18").Value = "27.360.41"
$ws.Range("E18").Value = "  -0.60%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "221.15"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.73%  "
$ws.Range("E20").Value = "  -0.07%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.628"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.02%  "
$ws.Range("E22").Value = "  -1.81%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.012"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.97%  "
$ws.Range("E24").Value = "  -0.11%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "144.91"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.48%  "
$ws.Range("E26").Value = "  -3.02%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1205"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.25%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.234"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.26%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "16.23"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.95%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05343"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.87%  "
$ws.Range("E31").Value = "  -1.03%  "
$ws.Range("E32").Value = "  -2.60%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.431"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.26%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.636"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.37%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.876"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.63%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9515"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.20%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.396"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.20%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5843"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.91%  "
$ws.Range("E39").Value = "  -1.02%  "
$ws.Range("D40").Value = "1.112.59"
$ws.Range("E40").Value = "  +5.30%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.799"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.95%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.004"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.14%  "
$ws.Range("E43").Value = "  -1.76%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "101.18"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.19%  "
$ws.Range("D45").Value = "1.854.34"
$ws.Range("E45").Value = "  -0.68%  "
$ws.Range("E46").Value = "  +1.91%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "57.59"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.20%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4543"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.31%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.005"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.06%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.118"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.44%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05225"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.42%  "
